$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd34"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 120.2294996666667
$ws.Range("H2").Value = 360.688499
$ws.Range("I2").Value = 0.2682304996487195
$ws.Range("J2").Value = 0.2682304996487195
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.450162
$ws.Range("N2").Value = 79.350486
$ws.Range("O2").Value = 0.9821351879331711
$ws.Range("P2").Value = 0.9821351879331711
$ws.Range("Q2").Value = 3180.08974336228
$ws.Range("R2").Value = 28620.80769026052
$ws.Range("S2").Value = 0.2634386121819036
$ws.Range("T2").Value = 0.2634386121819036

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd34"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 120.2294996666667
$ws.Range("H3").Value = 360.688499
$ws.Range("I3").Value = 0.2682304996487195
$ws.Range("J3").Value = 0.2682304996487195
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.321934
$ws.Range("N3").Value = 0.965802
$ws.Range("O3").Value = 0.01195390446349922
$ws.Range("P3").Value = 0.01195390446349922
$ws.Range("Q3").Value = 38.70596374568867
$ws.Range("R3").Value = 348.353673711198
$ws.Range("S3").Value = 0.003206401766997455
$ws.Range("T3").Value = 0.003206401766997455

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd34"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 120.2294996666667
$ws.Range("H4").Value = 360.688499
$ws.Range("I4").Value = 0.2682304996487195
$ws.Range("J4").Value = 0.2682304996487195
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1591883333333333
$ws.Range("N4").Value = 0.477565
$ws.Range("O4").Value = 0.005910907603329674
$ws.Range("P4").Value = 0.005910907603329674
$ws.Range("Q4").Value = 19.13913366943722
$ws.Range("R4").Value = 172.252203024935
$ws.Range("S4").Value = 0.001585485699818534
$ws.Range("T4").Value = 0.001585485699818534

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd34"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 320.4798073333333
$ws.Range("H5").Value = 961.4394219999999
$ws.Range("I5").Value = 0.7149864142051173
$ws.Range("J5").Value = 0.7149864142051173
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.450162
$ws.Range("N5").Value = 79.350486
$ws.Range("O5").Value = 0.9821351879331711
$ws.Range("P5").Value = 0.9821351879331711
$ws.Range("Q5").Value = 8476.742821695456
$ws.Range("R5").Value = 76290.68539525909
$ws.Range("S5").Value = 0.702213316285007
$ws.Range("T5").Value = 0.702213316285007

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cd34"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 320.4798073333333
$ws.Range("H6").Value = 961.4394219999999
$ws.Range("I6").Value = 0.7149864142051173
$ws.Range("J6").Value = 0.7149864142051173
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.321934
$ws.Range("N6").Value = 0.965802
$ws.Range("O6").Value = 0.01195390446349922
$ws.Range("P6").Value = 0.01195390446349922
$ws.Range("Q6").Value = 103.1733462940493
$ws.Range("R6").Value = 928.5601166464439
$ws.Range("S6").Value = 0.008546879288107856
$ws.Range("T6").Value = 0.008546879288107856

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cd34"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 320.4798073333333
$ws.Range("H7").Value = 961.4394219999999
$ws.Range("I7").Value = 0.7149864142051173
$ws.Range("J7").Value = 0.7149864142051173
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1591883333333333
$ws.Range("N7").Value = 0.477565
$ws.Range("O7").Value = 0.005910907603329674
$ws.Range("P7").Value = 0.005910907603329674
$ws.Range("Q7").Value = 51.01664639638111
$ws.Range("R7").Value = 459.14981756743
$ws.Range("S7").Value = 0.004226218632002448
$ws.Range("T7").Value = 0.004226218632002448

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd34"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.522716666666668
$ws.Range("H8").Value = 22.56815
$ws.Range("I8").Value = 0.01678308614616306
$ws.Range("J8").Value = 0.01678308614616306
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.450162
$ws.Range("N8").Value = 79.350486
$ws.Range("O8").Value = 0.9821351879331711
$ws.Range("P8").Value = 0.9821351879331711
$ws.Range("Q8").Value = 198.9770745134334
$ws.Range("R8").Value = 1790.7936706209
$ws.Range("S8").Value = 0.01648325946626046
$ws.Range("T8").Value = 0.01648325946626046

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd34"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.522716666666668
$ws.Range("H9").Value = 22.56815
$ws.Range("I9").Value = 0.01678308614616306
$ws.Range("J9").Value = 0.01678308614616306
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.321934
$ws.Range("N9").Value = 0.965802
$ws.Range("O9").Value = 0.01195390446349922
$ws.Range("P9").Value = 0.01195390446349922
$ws.Range("Q9").Value = 2.421818267366667
$ws.Range("R9").Value = 21.7963644063
$ws.Range("S9").Value = 0.0002006234083939106
$ws.Range("T9").Value = 0.0002006234083939106

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd34"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.522716666666668
$ws.Range("H10").Value = 22.56815
$ws.Range("I10").Value = 0.01678308614616306
$ws.Range("J10").Value = 0.01678308614616306
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1591883333333333
$ws.Range("N10").Value = 0.477565
$ws.Range("O10").Value = 0.005910907603329674
$ws.Range("P10").Value = 0.005910907603329674
$ws.Range("Q10").Value = 1.197528728305556
$ws.Range("R10").Value = 10.77775855475
$ws.Range("S10").Value = 0.00009920327150869218
$ws.Range("T10").Value = 0.00009920327150869218
